# Update structure of the input_data test fixture.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column ("values") used to be all 1s -> give each row a distinct value.
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 3
$ws.Range("D5").Value = 4
$ws.Range("D6").Value = 5

# Move the active selection to where the author left off editing.
$ws.Range("H11").Select() | Out-Null
